# Auto-generated Excel COM-interop script to apply Goblin_Profits sheet updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 53.666668  # ALC!H53: 65.625 -> 53.666668
$ws.Cells.Item(53, 9).Value = 53.666668  # ALC!I53: 57.363636 -> 53.666668
$ws.Cells.Item(53, 10).Value = 0  # ALC!J53: 83.8 -> 0
$ws.Cells.Item(53, 11).Value = 53.666668  # ALC!K53: 57.363636 -> 53.666668
$ws.Cells.Item(53, 12).Value = 0  # ALC!L53: 83.8 -> 0
$ws.Cells.Item(53, 13).Value = 583.333332  # ALC!M53: 579.636364 -> 583.333332
$ws.Cells.Item(53, 14).ClearContents()  # ALC!N53: -1357.8 -> (removed)

$ws.Cells.Item(137, 8).Value = 1268.6945  # ALC!H137: 1331.6471 -> 1268.6945
$ws.Cells.Item(137, 9).Value = 1249.9354  # ALC!I137: 1322.2069 -> 1249.9354
$ws.Cells.Item(137, 10).Value = 1385  # ALC!J137: 1386.4 -> 1385
$ws.Cells.Item(137, 11).Value = 3749.8062  # ALC!K137: 3966.620699999999 -> 3749.8062
$ws.Cells.Item(137, 12).Value = 4155  # ALC!L137: 4159.200000000001 -> 4155
$ws.Cells.Item(137, 13).Value = -1199.8062  # ALC!M137: -1416.620699999999 -> -1199.8062
$ws.Cells.Item(137, 14).Value = -9255  # ALC!N137: -9259.200000000001 -> -9255

$ws.Cells.Item(138, 8).Value = 1745.1111  # ALC!H138: 1723.4791 -> 1745.1111
$ws.Cells.Item(138, 10).Value = 1846.6765  # ALC!J138: 1810.3784 -> 1846.6765
$ws.Cells.Item(138, 12).Value = 5540.029500000001  # ALC!L138: 5431.135200000001 -> 5540.029500000001
$ws.Cells.Item(138, 14).Value = -15820.0295  # ALC!N138: -15711.1352 -> -15820.0295

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(41, 8).Value = 2433  # ARM!H41: 7126.25 -> 2433
$ws.Cells.Item(41, 9).Value = 2433  # ARM!I41: 7126.25 -> 2433
$ws.Cells.Item(41, 11).Value = 2433  # ARM!K41: 7126.25 -> 2433
$ws.Cells.Item(41, 13).Value = -2019  # ARM!M41: -6712.25 -> -2019

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2085.875  # BSM!H20: 2006.5834 -> 2085.875
$ws.Cells.Item(20, 9).Value = 1669.7142  # BSM!I20: 1625.8572 -> 1669.7142
$ws.Cells.Item(20, 10).Value = 4999  # BSM!J20: 2539.6 -> 4999
$ws.Cells.Item(20, 11).Value = 1669.7142  # BSM!K20: 1625.8572 -> 1669.7142
$ws.Cells.Item(20, 12).Value = 4999  # BSM!L20: 2539.6 -> 4999
$ws.Cells.Item(20, 13).Value = -1422.7142  # BSM!M20: -1378.8572 -> -1422.7142
$ws.Cells.Item(20, 14).Value = -5493  # BSM!N20: -3033.6 -> -5493

$ws.Cells.Item(36, 8).Value = 977.5  # BSM!H36: 1012.8333 -> 977.5
$ws.Cells.Item(36, 9).Value = 977.5  # BSM!I36: 1012.8333 -> 977.5
$ws.Cells.Item(36, 11).Value = 977.5  # BSM!K36: 1012.8333 -> 977.5
$ws.Cells.Item(36, 13).Value = -443.5  # BSM!M36: -478.8333 -> -443.5

$ws.Cells.Item(37, 8).Value = 3262.5  # BSM!H37: 1687.5 -> 3262.5
$ws.Cells.Item(37, 9).Value = 1683.3334  # BSM!I37: 2275 -> 1683.3334
$ws.Cells.Item(37, 10).Value = 8000  # BSM!J37: 1100 -> 8000
$ws.Cells.Item(37, 11).Value = 1683.3334  # BSM!K37: 2275 -> 1683.3334
$ws.Cells.Item(37, 12).Value = 8000  # BSM!L37: 1100 -> 8000
$ws.Cells.Item(37, 13).Value = -1546.3334  # BSM!M37: -2138 -> -1546.3334
$ws.Cells.Item(37, 14).Value = -8274  # BSM!N37: -1374 -> -8274

$ws.Cells.Item(49, 8).Value = 17099.8  # BSM!H49: 20499.75 -> 17099.8
$ws.Cells.Item(49, 9).Value = 3500  # BSM!I49: 0 -> 3500
$ws.Cells.Item(49, 11).Value = 3500  # BSM!K49: 0 -> 3500
$ws.Cells.Item(49, 13).Value = -3261  # BSM!M49: None -> -3261

$ws.Cells.Item(105, 8).Value = 3090.484  # BSM!H105: 3201.8333 -> 3090.484
$ws.Cells.Item(105, 9).Value = 2632.28  # BSM!I105: 2682.28 -> 2632.28
$ws.Cells.Item(105, 10).Value = 4999.6665  # BSM!J105: 5799.6 -> 4999.6665
$ws.Cells.Item(105, 11).Value = 2632.28  # BSM!K105: 2682.28 -> 2632.28
$ws.Cells.Item(105, 12).Value = 4999.6665  # BSM!L105: 5799.6 -> 4999.6665
$ws.Cells.Item(105, 13).Value = -885.2800000000002  # BSM!M105: -935.2800000000002 -> -885.2800000000002
$ws.Cells.Item(105, 14).Value = -8493.666499999999  # BSM!N105: -9293.6 -> -8493.666499999999

$ws.Cells.Item(134, 8).Value = 2393.4092  # BSM!H134: 2615.0527 -> 2393.4092
$ws.Cells.Item(134, 9).Value = 2332.75  # BSM!I134: 2569.7646 -> 2332.75
$ws.Cells.Item(134, 11).Value = 6998.25  # BSM!K134: 7709.293799999999 -> 6998.25
$ws.Cells.Item(134, 13).Value = -4463.25  # BSM!M134: -5174.293799999999 -> -4463.25

$ws.Cells.Item(137, 8).Value = 75000  # BSM!H137: 0 -> 75000
$ws.Cells.Item(137, 10).Value = 75000  # BSM!J137: 0 -> 75000
$ws.Cells.Item(137, 12).Value = 75000  # BSM!L137: 0 -> 75000
$ws.Cells.Item(137, 14).Value = -85200  # BSM!N137: None -> -85200

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(51, 8).Value = 26000  # CRP!H51: 12000 -> 26000
$ws.Cells.Item(51, 10).Value = 40000  # CRP!J51: 0 -> 40000
$ws.Cells.Item(51, 12).Value = 40000  # CRP!L51: 0 -> 40000
$ws.Cells.Item(51, 14).Value = -41472  # CRP!N51: None -> -41472

$ws.Cells.Item(61, 8).Value = 26000  # CRP!H61: 12000 -> 26000
$ws.Cells.Item(61, 10).Value = 40000  # CRP!J61: 0 -> 40000
$ws.Cells.Item(61, 12).Value = 40000  # CRP!L61: 0 -> 40000
$ws.Cells.Item(61, 14).Value = -40696  # CRP!N61: None -> -40696

$ws.Cells.Item(62, 8).Value = 16485  # CRP!H62: 14844.286 -> 16485
$ws.Cells.Item(62, 10).Value = 19251.25  # CRP!J62: 16401 -> 19251.25
$ws.Cells.Item(62, 12).Value = 19251.25  # CRP!L62: 16401 -> 19251.25
$ws.Cells.Item(62, 14).Value = -20499.25  # CRP!N62: -17649 -> -20499.25

$ws.Cells.Item(65, 8).Value = 16485  # CRP!H65: 14844.286 -> 16485
$ws.Cells.Item(65, 10).Value = 19251.25  # CRP!J65: 16401 -> 19251.25
$ws.Cells.Item(65, 12).Value = 96256.25  # CRP!L65: 82005 -> 96256.25
$ws.Cells.Item(65, 14).Value = -102496.25  # CRP!N65: -88245 -> -102496.25

$ws.Cells.Item(134, 8).Value = 3035.0527  # CRP!H134: 3050.4736 -> 3035.0527
$ws.Cells.Item(134, 9).Value = 2787.0557  # CRP!I134: 2803.3333 -> 2787.0557
$ws.Cells.Item(134, 11).Value = 8361.167099999999  # CRP!K134: 8409.999899999999 -> 8361.167099999999
$ws.Cells.Item(134, 13).Value = -5826.167099999999  # CRP!M134: -5874.999899999999 -> -5826.167099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 4562.25  # CUL!H18: 5142.5713 -> 4562.25
$ws.Cells.Item(18, 9).Value = 949.5  # CUL!I18: 1399 -> 949.5
$ws.Cells.Item(18, 11).Value = 2848.5  # CUL!K18: 4197 -> 2848.5
$ws.Cells.Item(18, 13).Value = -2679.5  # CUL!M18: -4028 -> -2679.5

$ws.Cells.Item(40, 8).Value = 75  # CUL!H40: 74 -> 75
$ws.Cells.Item(40, 9).Value = 50  # CUL!I40: 47.5 -> 50
$ws.Cells.Item(40, 10).Value = 87.5  # CUL!J40: 91.666664 -> 87.5
$ws.Cells.Item(40, 11).Value = 200  # CUL!K40: 190 -> 200
$ws.Cells.Item(40, 12).Value = 350  # CUL!L40: 366.666656 -> 350
$ws.Cells.Item(40, 13).Value = -131  # CUL!M40: -121 -> -131
$ws.Cells.Item(40, 14).Value = -488  # CUL!N40: -504.666656 -> -488

$ws.Cells.Item(51, 8).Value = 959.6  # CUL!H51: 1049.6666 -> 959.6

$ws.Cells.Item(64, 8).Value = 1266.3334  # CUL!H64: 1299.5 -> 1266.3334
$ws.Cells.Item(64, 9).Value = 1266.3334  # CUL!I64: 1299.5 -> 1266.3334
$ws.Cells.Item(64, 11).Value = 3799.0002  # CUL!K64: 3898.5 -> 3799.0002
$ws.Cells.Item(64, 13).Value = -3529.0002  # CUL!M64: -3628.5 -> -3529.0002

$ws.Cells.Item(67, 8).Value = 1266.3334  # CUL!H67: 1299.5 -> 1266.3334
$ws.Cells.Item(67, 9).Value = 1266.3334  # CUL!I67: 1299.5 -> 1266.3334
$ws.Cells.Item(67, 11).Value = 3799.0002  # CUL!K67: 3898.5 -> 3799.0002
$ws.Cells.Item(67, 13).Value = -2863.0002  # CUL!M67: -2962.5 -> -2863.0002

$ws.Cells.Item(99, 8).Value = 28170.666  # CUL!H99: 26856.715 -> 28170.666
$ws.Cells.Item(99, 9).Value = 20012.5  # CUL!I99: 40000 -> 20012.5
$ws.Cells.Item(99, 10).Value = 32249.75  # CUL!J99: 24666.166 -> 32249.75
$ws.Cells.Item(99, 11).Value = 60037.5  # CUL!K99: 120000 -> 60037.5
$ws.Cells.Item(99, 12).Value = 96749.25  # CUL!L99: 73998.49800000001 -> 96749.25
$ws.Cells.Item(99, 13).Value = -57791.5  # CUL!M99: -117754 -> -57791.5
$ws.Cells.Item(99, 14).Value = -101241.25  # CUL!N99: -78490.49800000001 -> -101241.25

$ws.Cells.Item(120, 8).Value = 34126.285  # CUL!H120: 37772.363 -> 34126.285
$ws.Cells.Item(120, 9).Value = 18642.166  # CUL!I120: 16799.2 -> 18642.166
$ws.Cells.Item(120, 10).Value = 45739.375  # CUL!J120: 55250 -> 45739.375
$ws.Cells.Item(120, 11).Value = 55926.49800000001  # CUL!K120: 50397.60000000001 -> 55926.49800000001
$ws.Cells.Item(120, 12).Value = 137218.125  # CUL!L120: 165750 -> 137218.125
$ws.Cells.Item(120, 13).Value = -51088.49800000001  # CUL!M120: -45559.60000000001 -> -51088.49800000001
$ws.Cells.Item(120, 14).Value = -146894.125  # CUL!N120: -175426 -> -146894.125

$ws.Cells.Item(125, 8).Value = 10000  # CUL!H125: 7500 -> 10000
$ws.Cells.Item(125, 9).Value = 10000  # CUL!I125: 7500 -> 10000
$ws.Cells.Item(125, 11).Value = 30000  # CUL!K125: 22500 -> 30000
$ws.Cells.Item(125, 13).Value = -25080  # CUL!M125: -17580 -> -25080

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value = 400  # GSM!H41: 199.5 -> 400
$ws.Cells.Item(41, 9).Value = 400  # GSM!I41: 199.5 -> 400
$ws.Cells.Item(41, 11).Value = 400  # GSM!K41: 199.5 -> 400
$ws.Cells.Item(41, 13).Value = -45  # GSM!M41: 155.5 -> -45

$ws.Cells.Item(70, 8).Value = 9675  # GSM!H70: 7599.4707 -> 9675
$ws.Cells.Item(70, 9).Value = 9566.666999999999  # GSM!I70: 7392.2856 -> 9566.666999999999
$ws.Cells.Item(70, 10).Value = 10000  # GSM!J70: 8566.333000000001 -> 10000
$ws.Cells.Item(70, 11).Value = 9566.666999999999  # GSM!K70: 7392.2856 -> 9566.666999999999
$ws.Cells.Item(70, 12).Value = 10000  # GSM!L70: 8566.333000000001 -> 10000
$ws.Cells.Item(70, 13).Value = -9296.666999999999  # GSM!M70: -7122.2856 -> -9296.666999999999
$ws.Cells.Item(70, 14).Value = -10540  # GSM!N70: -9106.333000000001 -> -10540

$ws.Cells.Item(73, 8).Value = 9675  # GSM!H73: 7599.4707 -> 9675
$ws.Cells.Item(73, 9).Value = 9566.666999999999  # GSM!I73: 7392.2856 -> 9566.666999999999
$ws.Cells.Item(73, 10).Value = 10000  # GSM!J73: 8566.333000000001 -> 10000
$ws.Cells.Item(73, 11).Value = 9566.666999999999  # GSM!K73: 7392.2856 -> 9566.666999999999
$ws.Cells.Item(73, 12).Value = 10000  # GSM!L73: 8566.333000000001 -> 10000
$ws.Cells.Item(73, 13).Value = -8630.666999999999  # GSM!M73: -6456.2856 -> -8630.666999999999
$ws.Cells.Item(73, 14).Value = -11872  # GSM!N73: -10438.333 -> -11872

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3722.963  # LTW!H22: 3827.7307 -> 3722.963
$ws.Cells.Item(22, 9).Value = 3572.1667  # LTW!I22: 3806.0908 -> 3572.1667
$ws.Cells.Item(22, 11).Value = 3572.1667  # LTW!K22: 3806.0908 -> 3572.1667
$ws.Cells.Item(22, 13).Value = -3277.1667  # LTW!M22: -3511.0908 -> -3277.1667

$ws.Cells.Item(27, 8).Value = 3722.963  # LTW!H27: 3827.7307 -> 3722.963
$ws.Cells.Item(27, 9).Value = 3572.1667  # LTW!I27: 3806.0908 -> 3572.1667
$ws.Cells.Item(27, 11).Value = 3572.1667  # LTW!K27: 3806.0908 -> 3572.1667
$ws.Cells.Item(27, 13).Value = -3465.1667  # LTW!M27: -3699.0908 -> -3465.1667

$ws.Cells.Item(122, 8).Value = 4036.7646  # LTW!H122: 4164.0625 -> 4036.7646
$ws.Cells.Item(122, 9).Value = 4862.6  # LTW!I122: 5180.6665 -> 4862.6
$ws.Cells.Item(122, 11).Value = 14587.8  # LTW!K122: 15541.9995 -> 14587.8
$ws.Cells.Item(122, 13).Value = -12137.8  # LTW!M122: -13091.9995 -> -12137.8

$ws.Cells.Item(136, 8).Value = 13121.6455  # LTW!H136: 13631.348 -> 13121.6455
$ws.Cells.Item(136, 9).Value = 1821.1111  # LTW!I136: 1941.8572 -> 1821.1111
$ws.Cells.Item(136, 11).Value = 5463.3333  # LTW!K136: 5825.571599999999 -> 5463.3333
$ws.Cells.Item(136, 13).Value = -2913.3333  # LTW!M136: -3275.571599999999 -> -2913.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(33, 8).Value = 21402.715  # WVR!H33: 21999.6 -> 21402.715
$ws.Cells.Item(33, 10).Value = 21402.715  # WVR!J33: 21999.6 -> 21402.715
$ws.Cells.Item(33, 12).Value = 21402.715  # WVR!L33: 21999.6 -> 21402.715
$ws.Cells.Item(33, 14).Value = -21902.715  # WVR!N33: -22499.6 -> -21902.715

$ws.Cells.Item(36, 8).Value = 21402.715  # WVR!H36: 21999.6 -> 21402.715
$ws.Cells.Item(36, 10).Value = 21402.715  # WVR!J36: 21999.6 -> 21402.715
$ws.Cells.Item(36, 12).Value = 21402.715  # WVR!L36: 21999.6 -> 21402.715
$ws.Cells.Item(36, 14).Value = -21902.715  # WVR!N36: -22499.6 -> -21902.715

$ws.Cells.Item(64, 8).Value = 185000  # WVR!H64: 184999.5 -> 185000
$ws.Cells.Item(64, 10).Value = 185000  # WVR!J64: 184999.5 -> 185000
$ws.Cells.Item(64, 12).Value = 185000  # WVR!L64: 184999.5 -> 185000
$ws.Cells.Item(64, 14).Value = -185496  # WVR!N64: -185495.5 -> -185496

$ws.Cells.Item(67, 8).Value = 185000  # WVR!H67: 184999.5 -> 185000
$ws.Cells.Item(67, 10).Value = 185000  # WVR!J67: 184999.5 -> 185000
$ws.Cells.Item(67, 12).Value = 185000  # WVR!L67: 184999.5 -> 185000
$ws.Cells.Item(67, 14).Value = -186716  # WVR!N67: -186715.5 -> -186716

$ws.Cells.Item(136, 8).Value = 1002.6  # WVR!H136: 1009.2 -> 1002.6
$ws.Cells.Item(136, 9).Value = 1014.0526  # WVR!I136: 1021 -> 1014.0526
$ws.Cells.Item(136, 11).Value = 3042.1578  # WVR!K136: 3063 -> 3042.1578
$ws.Cells.Item(136, 13).Value = -492.1578  # WVR!M136: -513 -> -492.1578
